# Updated Nithin as Director of Fundraising
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$found = $ws.Cells.Find("Winnie Qi")
if ($found -ne $null) {
    $found.Value = "Nithin Senthil"
}
